{"js": "// Replace each two-digit multiplication expression with its updated value.\n// Each \"AxB=\" string in this document is unique, so an exact-match search\n// safely targets the single corresponding run in the table.\nconst replacements = [\n  [\"76\u00d723=\", \"13\u00d712=\"],\n  [\"46\u00d755=\", \"39\u00d767=\"],\n  [\"20\u00d776=\", \"98\u00d746=\"],\n  [\"76\u00d777=\", \"62\u00d779=\"],\n  [\"76\u00d726=\", \"17\u00d748=\"],\n  [\"36\u00d774=\", \"49\u00d721=\"],\n  [\"79\u00d740=\", \"20\u00d751=\"],\n  [\"13\u00d777=\", \"59\u00d749=\"],\n  [\"25\u00d756=\", \"36\u00d764=\"],\n  [\"12\u00d747=\", \"96\u00d789=\"],\n  [\"83\u00d786=\", \"77\u00d729=\"],\n  [\"18\u00d760=\", \"68\u00d786=\"],\n  [\"27\u00d766=\", \"99\u00d789=\"],\n  [\"27\u00d738=\", \"94\u00d722=\"],\n  [\"42\u00d731=\", \"32\u00d793=\"],\n  [\"33\u00d772=\", \"88\u00d716=\"],\n  [\"70\u00d798=\", \"23\u00d782=\"],\n  [\"39\u00d790=\", \"26\u00d749=\"],\n  [\"64\u00d778=\", \"23\u00d795=\"],\n  [\"68\u00d717=\", \"19\u00d772=\"],\n  [\"64\u00d772=\", \"34\u00d754=\"],\n  [\"85\u00d783=\", \"52\u00d780=\"],\n  [\"90\u00d763=\", \"36\u00d723=\"],\n  [\"51\u00d763=\", \"14\u00d724=\"],\n  [\"98\u00d754=\", \"59\u00d782=\"]\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Could not find text to replace: \" + oldText);\n  }\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update each two-digit multiplication expression to its new value.\n# wdReplaceAll (=2) ensures every matching occurrence in the body is updated;\n# each left-hand string below is unique in this document.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"76\u00d723=\"; New = \"13\u00d712=\" }\n    @{ Old = \"46\u00d755=\"; New = \"39\u00d767=\" }\n    @{ Old = \"20\u00d776=\"; New = \"98\u00d746=\" }\n    @{ Old = \"76\u00d777=\"; New = \"62\u00d779=\" }\n    @{ Old = \"76\u00d726=\"; New = \"17\u00d748=\" }\n    @{ Old = \"36\u00d774=\"; New = \"49\u00d721=\" }\n    @{ Old = \"79\u00d740=\"; New = \"20\u00d751=\" }\n    @{ Old = \"13\u00d777=\"; New = \"59\u00d749=\" }\n    @{ Old = \"25\u00d756=\"; New = \"36\u00d764=\" }\n    @{ Old = \"12\u00d747=\"; New = \"96\u00d789=\" }\n    @{ Old = \"83\u00d786=\"; New = \"77\u00d729=\" }\n    @{ Old = \"18\u00d760=\"; New = \"68\u00d786=\" }\n    @{ Old = \"27\u00d766=\"; New = \"99\u00d789=\" }\n    @{ Old = \"27\u00d738=\"; New = \"94\u00d722=\" }\n    @{ Old = \"42\u00d731=\"; New = \"32\u00d793=\" }\n    @{ Old = \"33\u00d772=\"; New = \"88\u00d716=\" }\n    @{ Old = \"70\u00d798=\"; New = \"23\u00d782=\" }\n    @{ Old = \"39\u00d790=\"; New = \"26\u00d749=\" }\n    @{ Old = \"64\u00d778=\"; New = \"23\u00d795=\" }\n    @{ Old = \"68\u00d717=\"; New = \"19\u00d772=\" }\n    @{ Old = \"64\u00d772=\"; New = \"34\u00d754=\" }\n    @{ Old = \"85\u00d783=\"; New = \"52\u00d780=\" }\n    @{ Old = \"90\u00d763=\"; New = \"36\u00d723=\" }\n    @{ Old = \"51\u00d763=\"; New = \"14\u00d724=\" }\n    @{ Old = \"98\u00d754=\"; New = \"59\u00d782=\" }\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair.Old\n    $find.Replacement.Text = $pair.New\n    $find.Execute(\n        $find.Text,\n        $false,\n        $false,\n        $false,\n        $false,\n        $false,\n        $true,\n        1,\n        $false,\n        $find.Replacement.Text,\n        2\n    ) | Out-Null\n}\n"}
